$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.884.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.07%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.811.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.70%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.49%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'309.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.54%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -0.37%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4940"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.63%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3885"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.18%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.09757"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +25.82%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.63%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'40.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.15%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'6.428"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.10%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "'Solana"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'20.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.64%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "'BinanceUSD"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'1.001"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.41%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'1.807.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.96%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'7.293"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.69%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +6.63%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'92.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.80%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06601"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.54%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.34%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'17.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.72%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.939"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.99%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'27.934.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.06%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.16%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.27%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'157.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.94%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'20.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.01%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'2.015.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.66%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.393"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.20%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'127.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.29%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.1060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.58%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.71%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.578"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.15%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.632"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.06753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.66%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'9.018"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.96%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.02321"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.88%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.2130"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.50%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'4.939"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.17%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'11.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.46%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.6183"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.77%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.9999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.35%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.145"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.21%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'13.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.52%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'Decentraland"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.5866"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.55%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'PancakeSwap"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'3.694"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.77%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.282"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.71%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'122.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.06%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.935"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.65%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.175"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.28%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.06792"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.29%  "
$ws.Range("E51").Style = "Normal"
